$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so numeric-looking strings
# (e.g. "27.892.98", "87.27") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.892.98"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").Value = "1.879.52"
$ws.Range("E3").Value = "  +1.59%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "333.05"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "0.4737"
$ws.Range("E7").Value = "  +5.42%  "

$ws.Range("E8").Value = "  +3.73%  "

$ws.Range("D9").Value = "48.16"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").Value = "0.08062"
$ws.Range("E10").Value = "  +2.72%  "

$ws.Range("D12").Value = "21.88"
$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("D13").Value = "1.897.95"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").Value = "7.211"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "87.27"
$ws.Range("E17").Value = "  +2.03%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001049"
$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("D19").Value = "0.06625"
$ws.Range("E19").Value = "  +2.03%  "

$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "27.990.88"
$ws.Range("E22").Value = "  +2.46%  "

$ws.Range("D23").Value = "5.520"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  +2.57%  "

$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  +2.33%  "

$ws.Range("D26").Value = "2.126.45"
$ws.Range("E26").Value = "  +2.94%  "

$ws.Range("D27").Value = "157.81"
$ws.Range("E27").Value = "  +4.35%  "

$ws.Range("D28").Value = "20.27"
$ws.Range("E28").Value = "  +4.74%  "

$ws.Range("E29").Value = "  +3.00%  "

$ws.Range("D30").Value = "5.639"
$ws.Range("E30").Value = "  +2.02%  "

$ws.Range("D31").Value = "122.66"
$ws.Range("E31").Value = "  +2.41%  "

$ws.Range("D32").Value = "0.9913"
$ws.Range("E32").Value = "  +6.33%  "

$ws.Range("E33").Value = "  +2.90%  "

$ws.Range("D34").Value = "1.466"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("D35").Value = "3.611"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("E37").Value = "  +2.52%  "

$ws.Range("D38").Value = "0.06127"
$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("D39").Value = "1.238"
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("D40").Value = "8.263"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").Value = "0.6034"
$ws.Range("E41").Value = "  +2.37%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "0.1912"
$ws.Range("E43").Value = "  +3.45%  "

$ws.Range("D44").Value = "10.36"

$ws.Range("D45").Value = "1.275"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("D46").Value = "0.5723"
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "12.31"
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("D48").Value = "1.954"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").Value = "3.409"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "0.06837"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").Value = "113.79"
$ws.Range("E51").Value = "  +5.27%  "
